# Apply scheduled-runner price/profit updates to Lich_Profits workbook
$wb = $excel.ActiveWorkbook

# ALC!row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1592.7273
$ws.Range("I11").Value = 1592.7273
$ws.Range("K11").Value = 1592.7273
$ws.Range("M11").Value = -1452.7273

# ALC!row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1577.4584
$ws.Range("J19").Value = 2005.8572
$ws.Range("L19").Value = 2005.8572
$ws.Range("N19").Value = -2355.8572

# ALC!row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2775.25
$ws.Range("I32").Value = 2950.5
$ws.Range("J32").Value = 2600
$ws.Range("K32").Value = 2950.5
$ws.Range("L32").Value = 2600
$ws.Range("M32").Value = -2624.5
$ws.Range("N32").Value = -3252

# ALC!row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 10424334
$ws.Range("I62").Value = 27785108
$ws.Range("J62").Value = 7868.7
$ws.Range("K62").Value = 27785108
$ws.Range("L62").Value = 7868.7
$ws.Range("M62").Value = -27784484
$ws.Range("N62").Value = -9116.700000000001

# ALC!row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3948
$ws.Range("I64").Value = 3597.3333
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 3597.3333
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -3349.3333
$ws.Range("N64").Value = -5496

# ALC!row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 10424334
$ws.Range("I65").Value = 27785108
$ws.Range("J65").Value = 7868.7
$ws.Range("K65").Value = 138925540
$ws.Range("L65").Value = 39343.5
$ws.Range("M65").Value = -138922420
$ws.Range("N65").Value = -45583.5

# ALC!row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3948
$ws.Range("I67").Value = 3597.3333
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 3597.3333
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -2739.3333
$ws.Range("N67").Value = -6716

# ALC!row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# ALC!row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2362.756
$ws.Range("I132").Value = 2299.35
$ws.Range("K132").Value = 6898.049999999999
$ws.Range("M132").Value = -4368.049999999999

# ALC!row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1132.0588
$ws.Range("I135").Value = 1200.5625
$ws.Range("J135").Value = 36
$ws.Range("K135").Value = 10805.0625
$ws.Range("L135").Value = 324
$ws.Range("M135").Value = -8270.0625
$ws.Range("N135").Value = -5394

# ARM!row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 929.6129
$ws.Range("I102").Value = 859.9655
$ws.Range("K102").Value = 859.9655
$ws.Range("M102").Value = 762.0345

# ARM!row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3550.2
$ws.Range("J122").Value = 3599.2
$ws.Range("L122").Value = 10797.6
$ws.Range("N122").Value = -15697.6

# BSM!row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1843.5834
$ws.Range("I107").Value = 1843.5834
$ws.Range("K107").Value = 1843.5834
$ws.Range("M107").Value = 76.41660000000002

# CRP!row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2604.2888
$ws.Range("I134").Value = 2573.4146
$ws.Range("K134").Value = 7720.2438
$ws.Range("M134").Value = -5185.2438

# CUL!row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 126.542854
$ws.Range("J2").Value = 147.36
$ws.Range("L2").Value = 884.1600000000001
$ws.Range("N2").Value = -1110.16

# CUL!row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 143868.58
$ws.Range("I7").Value = 250147.5
$ws.Range("J7").Value = 2163.3333
$ws.Range("K7").Value = 750442.5
$ws.Range("L7").Value = 6489.999899999999
$ws.Range("M7").Value = -750330.5
$ws.Range("N7").Value = -6713.999899999999

# GSM!row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3080.5833
$ws.Range("I70").Value = 2957.8
$ws.Range("J70").Value = 3694.5
$ws.Range("K70").Value = 2957.8
$ws.Range("L70").Value = 3694.5
$ws.Range("M70").Value = -2687.8
$ws.Range("N70").Value = -4234.5

# GSM!row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 3080.5833
$ws.Range("I73").Value = 2957.8
$ws.Range("J73").Value = 3694.5
$ws.Range("K73").Value = 2957.8
$ws.Range("L73").Value = 3694.5
$ws.Range("M73").Value = -2021.8
$ws.Range("N73").Value = -5566.5

# GSM!row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4370.879
$ws.Range("J80").Value = 4331.136
$ws.Range("L80").Value = 4331.136
$ws.Range("N80").Value = -6327.136

# GSM!row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4370.879
$ws.Range("J83").Value = 4331.136
$ws.Range("L83").Value = 21655.68
$ws.Range("N83").Value = -31639.68

# LTW!row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5494.778
$ws.Range("I7").Value = 4923.5713
$ws.Range("J7").Value = 7494
$ws.Range("K7").Value = 4923.5713
$ws.Range("L7").Value = 7494
$ws.Range("M7").Value = -4811.5713
$ws.Range("N7").Value = -7718

# LTW!row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15344.042
$ws.Range("I40").Value = 17463.947
$ws.Range("J40").Value = 7288.4
$ws.Range("K40").Value = 17463.947
$ws.Range("L40").Value = 7288.4
$ws.Range("M40").Value = -17327.947
$ws.Range("N40").Value = -7560.4

# LTW!row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5494.778
$ws.Range("I126").Value = 4923.5713
$ws.Range("J126").Value = 7494
$ws.Range("K126").Value = 14770.7139
$ws.Range("L126").Value = 22482
$ws.Range("M126").Value = -12300.7139
$ws.Range("N126").Value = -27422

# LTW!row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3016.0667
$ws.Range("I136").Value = 2964.9312
$ws.Range("J136").Value = 4499
$ws.Range("K136").Value = 8894.793600000001
$ws.Range("L136").Value = 13497
$ws.Range("M136").Value = -6344.793600000001
$ws.Range("N136").Value = -18597

# WVR!row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 3000
$ws.Range("I54").Value = 3000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 3000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -2480
$ws.Range("N54").ClearContents()

# WVR!row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1426.7778
$ws.Range("I100").Value = 1467.25
$ws.Range("K100").Value = 2934.5
$ws.Range("M100").Value = -2393.5

# WVR!row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1039.9642
$ws.Range("I113").Value = 887.7727
$ws.Range("K113").Value = 2663.3181
$ws.Range("M113").Value = -493.3181

# WVR!row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2689.3
$ws.Range("J122").Value = 2397.3635
$ws.Range("L122").Value = 7192.0905
$ws.Range("N122").Value = -12092.0905

# WVR!row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2361
$ws.Range("I126").Value = 2241.9375
$ws.Range("J126").Value = 2678.5
$ws.Range("K126").Value = 6725.8125
$ws.Range("L126").Value = 8035.5
$ws.Range("M126").Value = -4255.8125
$ws.Range("N126").Value = -12975.5
